# Refresh the "cryptos" price list with the latest scraped snapshot
# (GitHub Actions job, "Updated cryptos list ... with GitHub Actions").
# Most edits are plain Price (D) / Volume(1h) (E) text updates; a couple
# of D cells hold numeric-looking text (e.g. "0.999", "555.55") so we
# force the Text format on those specific cells before assigning, to
# keep them stored as text instead of being auto-coerced to numbers.
# Rows 29/30 additionally swap which coin (Bittensor vs
# Binance-PegBSC-USD) occupies which rank, so B/C/D/E are all rewritten
# for those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.846.77"
$ws.Range("E2").Value = "  -2.94%  "
$ws.Range("D3").Value = "2.492.61"
$ws.Range("E3").Value = "  -5.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.55"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.43"
$ws.Range("E6").Value = "  -4.98%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").Value = "2.492.13"
$ws.Range("E9").Value = "  -5.22%  "
$ws.Range("E10").Value = "  -7.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.46"
$ws.Range("E14").Value = "  -6.76%  "
$ws.Range("D15").Value = "2.942.10"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("E16").Value = "  -7.36%  "
$ws.Range("D17").Value = "61.771.70"
$ws.Range("D18").Value = "2.496.75"
$ws.Range("E18").Value = "  -4.98%  "
$ws.Range("E19").Value = "  -7.37%  "
$ws.Range("E20").Value = "  -7.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  -6.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "324.03"
$ws.Range("E22").Value = "  -6.06%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -4.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.27"
$ws.Range("E25").Value = "  -5.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000100"
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("E27").Value = "  -3.89%  "
$ws.Range("D28").Value = "2.611.30"
$ws.Range("E28").Value = "  -5.02%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "539.22"
$ws.Range("E30").Value = "  -11.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  -8.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.62"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E33").Value = "  -5.29%  "
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("E35").Value = "  -7.95%  "
$ws.Range("E36").Value = "  -9.38%  "
$ws.Range("E37").Value = "  -8.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("E40").Value = "  -5.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.66"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.43"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("E45").Value = "  -6.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.44"
$ws.Range("E46").Value = "  -6.18%  "
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.21"
$ws.Range("E48").Value = "  -14.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0541"
$ws.Range("E49").Value = "  -7.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.600"
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("E51").Value = "  -4.71%  "
